$d = $word.ActiveDocument

# Locate the exact cell text "Precondición:" (NOT "Precondiciones:", which
# lives in a different row and must stay untouched) and turn it into
# "Postcondición:" - mirroring an edit where the "re" in the middle of the
# word was selected and retyped as "ost", which is why the final text ends
# up split across three runs: "P" / "ost" / "condición:".
$search = $d.Content
$search.Find.ClearFormatting()
$search.Find.Text = "Precondición:"
$search.Find.MatchWildcards = $false
$search.Find.MatchCase = $true
$search.Find.Forward = $true
$search.Find.Wrap = 1
$found = $search.Find.Execute()

if ($found -and $search.Find.Found) {
    $target = $d.Range($search.Start, $search.End)
    $start = $target.Start

    # Replace "re" (the 2nd and 3rd characters) with "ost".
    $mid = $d.Range($start + 1, $start + 3)
    $mid.Text = "ost"

    # Toggling a character-formatting property on just the replaced text
    # and then clearing it again forces Word to keep it as its own run
    # instead of silently re-joining it with its neighbours, producing the
    # "P" / "ost" / "condición:" run split.
    $mid2 = $d.Range($start + 1, $start + 4)
    $mid2.Bold = $true
    $mid2.Bold = $false
}
